$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.340.03'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '3.611.88'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('E4').Value = '  -0.06%  '
$cell = $ws.Range('D5')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '604.55'
$cell.Style = $savedStyle
$ws.Range('E5').Value = '  +0.23%  '
$cell = $ws.Range('D6')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '196.10'
$cell.Style = $savedStyle
$ws.Range('E6').Value = '  -0.74%  '
$cell = $ws.Range('D7')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.626'
$cell.Style = $savedStyle
$ws.Range('E7').Value = '  -0.26%  '
$cell = $ws.Range('D8')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = $savedStyle
$ws.Range('E9').Value = '  -1.55%  '
$cell = $ws.Range('D10')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.647'
$cell.Style = $savedStyle
$ws.Range('E10').Value = '  -1.22%  '
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('E13').Value = '  -0.19%  '
$ws.Range('D14').Value = '4.181.59'
$ws.Range('E14').Value = '  +1.55%  '
$cell = $ws.Range('D15')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '13.02'
$cell.Style = $savedStyle
$ws.Range('E15').Value = '  +2.46%  '
$cell = $ws.Range('D16')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '596.32'
$cell.Style = $savedStyle
$ws.Range('E16').Value = '  -1.57%  '
$ws.Range('D17').Value = '70.410.25'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Range('D18')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '19.04'
$cell.Style = $savedStyle
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.591.78'
$ws.Range('E19').Value = '  +0.86%  '
$ws.Range('E20').Value = '  +1.41%  '
$cell = $ws.Range('D21')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.997'
$cell.Style = $savedStyle
$ws.Range('E21').Value = '  -0.19%  '
$cell = $ws.Range('D22')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '17.86'
$cell.Style = $savedStyle
$ws.Range('E22').Value = '  -0.95%  '
$cell = $ws.Range('D23')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '5.19'
$cell.Style = $savedStyle
$ws.Range('E23').Value = '  -0.74%  '
$cell = $ws.Range('D24')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '101.85'
$cell.Style = $savedStyle
$ws.Range('E24').Value = '  -1.31%  '
$cell = $ws.Range('D25')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '4.62'
$cell.Style = $savedStyle
$ws.Range('E25').Value = '  +0.16%  '
$ws.Range('E26').Value = '  -3.18%  '
$cell = $ws.Range('D27')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '10.77'
$cell.Style = $savedStyle
$ws.Range('E27').Value = '  -1.80%  '
$cell = $ws.Range('D28')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '9.60'
$cell.Style = $savedStyle
$ws.Range('E28').Value = '  -0.64%  '
$cell = $ws.Range('D29')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '33.79'
$cell.Style = $savedStyle
$ws.Range('E29').Value = '  -0.03%  '
$cell = $ws.Range('D30')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '4.75'
$cell.Style = $savedStyle
$ws.Range('E30').Value = '  +5.68%  '
$cell = $ws.Range('D31')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '7.23'
$cell.Style = $savedStyle
$ws.Range('E31').Value = '  +0.89%  '
$ws.Range('E32').Value = '  -3.23%  '
$cell = $ws.Range('D33')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.117'
$cell.Style = $savedStyle
$ws.Range('E33').Value = '  +1.47%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$cell = $ws.Range('D34')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '63.19'
$cell.Style = $savedStyle
$ws.Range('E34').Value = '  -0.43%  '
$ws.Range('B35').Value = 'PEPE'
$ws.Range('C35').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D35').Value = '0.0₃0897'
$ws.Range('E35').Value = '  +6.46%  '
$ws.Range('D36').Value = '3.895.87'
$ws.Range('E36').Value = '  +3.32%  '
$cell = $ws.Range('D37')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '543.80'
$cell.Style = $savedStyle
$ws.Range('E37').Value = '  +11.96%  '
$cell = $ws.Range('D38')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.14'
$cell.Style = $savedStyle
$ws.Range('E38').Value = '  +1.55%  '
$ws.Range('E39').Value = '  +0.00%  '
$cell = $ws.Range('D40')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '36.88'
$cell.Style = $savedStyle
$ws.Range('E40').Value = '  -0.03%  '
$cell = $ws.Range('D41')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.391'
$cell.Style = $savedStyle
$ws.Range('E41').Value = '  -1.30%  '
$ws.Range('E42').Value = '  -3.72%  '
$ws.Range('E43').Value = '  -1.54%  '
$ws.Range('E44').Value = '  -0.59%  '
$cell = $ws.Range('D45')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.43'
$cell.Style = $savedStyle
$ws.Range('E45').Value = '  +3.07%  '
$cell = $ws.Range('D46')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.87'
$cell.Style = $savedStyle
$ws.Range('E46').Value = '  +0.45%  '
$ws.Range('E47').Value = '  -0.28%  '
$cell = $ws.Range('D48')
$savedStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '8.61'
$cell.Style = $savedStyle
$ws.Range('E48').Value = '  -0.97%  '
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('E50').Value = '  -0.74%  '
$ws.Range('E51').Value = '  +0.02%  '
